$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table occupies rows 16-23 (columns B:G hold values; H:J are blank
# style-only cells). The edit reverses the order of these 8 records so the
# "GERARDO RAFAEL ACUÑA GONZALEZ" periods (2305..2211) now come first and the
# "RONALD ALBERTO MARQUEZ BARRAZA" record (1912) comes last.

$firstRow = 16
$lastRow = 23

# Snapshot current values for columns B through G across the block.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 2; $c -le 7; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write back in reverse row order.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $lastRow - ($r - $firstRow)
    $rowVals = $snapshot[$srcRow]
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $rowVals[$c]
    }
}
